{"js": "// Replace the hyperlinked YouTube URL with a plain-text (non-hyperlinked) URL.\nconst oldUrl = \"https://youtu.be/BF_oIouqygQ\";\nconst newUrl = \"https://youtu.be/zT_skd-rYSg\";\n\nconst results = context.document.body.search(oldUrl, { matchCase: true, matchWildcards: false });\nresults.load(\"text\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error(\"Could not find the hyperlink text to replace: \" + oldUrl);\n}\n\nconst target = results.items[0];\n\n// Clearing the hyperlink also drops the Hyperlink character style / formatting,\n// turning the run back into plain text.\ntarget.hyperlink = \"\";\n\n// Swap the old URL text for the new one in-place.\ntarget.insertText(newUrl, \"Replace\");\n\nawait context.sync();\n", "ps1": "# Replace the hyperlinked YouTube URL with a plain-text (non-hyperlinked) URL.\n$d = $word.ActiveDocument\n$oldUrl = \"https://youtu.be/BF_oIouqygQ\"\n$newUrl = \"https://youtu.be/zT_skd-rYSg\"\n\n# Remove the hyperlink field itself (keeps its display text in place for now).\nforeach ($h in $d.Hyperlinks) {\n    if ($h.Address -eq $oldUrl) {\n        $h.Delete()\n        break\n    }\n}\n\n# Find the (now unlinked) URL text and swap it out for fresh, unstyled text so\n# the old \"Hyperlink\" character style doesn't carry over to the replacement.\n$range = $d.Content\n$found = $range.Find.Execute($oldUrl)\nif ($found) {\n    $range.Delete()\n    $range.InsertAfter($newUrl)\n}\n"}
